$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arduino_avr_uno")

# --- Add the new row of results for release v6.16 ---
$ws.Range("A33").Value = "v6.16"
$ws.Range("B33").Value = 44044
$ws.Range("C33").Value = 6930
$ws.Range("D33").Value = 5738
$ws.Range("E33").Value = 40

# Copy the date cell's format from the row above so B33 keeps the same
# date style (instead of getting a brand-new style entry).
$ws.Range("B32").Copy()
$ws.Range("B33").PasteSpecial(-4122)

# --- Extend the "Code size on ATmega328" chart's series ranges so the
# new row is included in the plot (JsonParserExample / JsonGeneratorExample
# series, columns C and D). ---
$cos = $ws.ChartObjects()
$chart1 = $cos.Item(1).Chart
$series1 = $chart1.SeriesCollection()

$s1 = $series1.Item(1)
$s1.Formula = "=SERIES(arduino_avr_uno!`$C`$1,arduino_avr_uno!`$A`$2:`$A`$33,arduino_avr_uno!`$C`$2:`$C`$33,1)"

$s2 = $series1.Item(2)
$s2.Formula = "=SERIES(arduino_avr_uno!`$D`$1,arduino_avr_uno!`$A`$2:`$A`$33,arduino_avr_uno!`$D`$2:`$D`$33,2)"

# --- Reflect the author's resulting UI state: scrolled down, zoomed to
# 115%, with C31 selected. ---
$ws.Range("C31").Select()
$excel.ActiveWindow.Zoom = 115
